$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "2024-06-14 17:28:15"
$ws.Range("D12").Value = 200
$ws.Range("E12").Value = 6

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = "2024-06-14 17:28:16"
$ws.Range("D13").Value = 200
$ws.Range("E13").Value = 6
